$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / URL / label updates (safe as-is, Excel will not reinterpret as numbers) ---
$ws.Range("D2").Value = "70.086.82"
$ws.Range("E2").Value = "  +2.67%  "
$ws.Range("D3").Value = "2.455.83"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +12.40%  "
$ws.Range("D10").Value = "2.447.52"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("E14").Value = "  +8.33%  "
$ws.Range("D15").Value = "69.953.88"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "2.900.19"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  +5.18%  "
$ws.Range("D18").Value = "2.450.00"
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("B19").Value = "Binance-PegBSC-USD"
$ws.Range("C19").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("E19").Value = "  +69.12%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E20").Value = "  +5.23%  "
$ws.Range("E21").Value = "  +5.74%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E22").Value = "  +2.12%  "
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("E24").Value = "  +6.25%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E27").Value = "  +6.55%  "
$ws.Range("D28").Value = "2.573.71"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E29").Value = "  +5.61%  "
$ws.Range("D30").Value = "0.0₃0870"
$ws.Range("E30").Value = "  +7.85%  "
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("E32").Value = "  +10.38%  "
$ws.Range("E33").Value = "  +6.23%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("E42").Value = "  +6.18%  "
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("E45").Value = "  +8.59%  "
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("E49").Value = "  +3.47%  "
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("E51").Value = "  +1.80%  "

# --- Numeric-looking text values: force as Text so they stay stored as strings (matches source formatting) ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.512"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.176"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000183"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "343.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "450.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.110"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.305"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0724"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.491"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.563"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0933"
$ws.Range("D51").Style = "Normal"
